$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced
# to stay text (matching the source inlineStr cells) by temporarily
# applying a text number format, same trick real Excel users rely on
# (leading apostrophe / Format Cells > Text) to stop auto-number coercion.
# ClearFormats() afterwards drops the temporary "@" format again -- safe
# here because none of these cells carry any pre-existing styling.

$ws.Range("D2").Value = '28.020.43'
$ws.Range("E2").Value = '  +1.45%  '

$ws.Range("D3").Value = '1.642.57'
$ws.Range("E3").Value = '  +0.56%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.524'
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.69'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.64%  '

$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0615'
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.93%  '

$ws.Range("D12").Value = '1.875.80'
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("D13").Value = '1.639.01'
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.45%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.577'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.03'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("D17").Value = '28.004.53'
$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.18'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.64'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.75'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.56'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.98'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.76'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.52%  '

$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("E30").Value = '  +0.87%  '

$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("E32").Value = '  +2.19%  '

$ws.Range("E33").Value = '  +0.39%  '

$ws.Range("D34").Value = '1.406.64'
$ws.Range("E34").Value = '  -4.67%  '

$ws.Range("E35").Value = '  +2.62%  '

$ws.Range("E36").Value = '  +0.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.892'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0169'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.557'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.917'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.53%  '

$ws.Range("E41").Value = '  -0.75%  '

$ws.Range("E42").Value = '  -0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.86'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +7.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.24%  '

$ws.Range("E45").Value = '  +2.37%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").Value = '1.784.34'
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.15'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.78%  '
